# Edit accounts workbook: append newly tracked advertiser accounts and
# renumber the handoff of the existing "가이드맨" row as it shifts down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new account ("키들") was inserted right before the previous last data
# row, pushing the old row 20 ("가이드맨") down to row 21 — use a real
# row insert so existing data shifts instead of being overwritten.
$ws.Rows.Item(20).Insert()

# rows keyed as: row number, company name (A), ad spend amount (B), manager name (C)
$newRows = @(
    @(20, "키들", 3322328, "승훈"),
    @(22, "청주오스코", 3300114, "미혜"),
    @(23, "리틀코리아", 207750, "미혜"),
    @(24, "부성에이티", 3113312, "미혜"),
    @(25, "컨텍스쳐", 2404793, "미혜"),
    @(26, "오르비", 2026164, "미혜"),
    @(27, "아틀라시안", 1340565, "미혜"),
    @(28, "이유즈", 2760916, "정인"),
    @(29, "수협보험", 3217504, "정인"),
    @(30, "네모조명", 2378803, "정인"),
    @(31, "아뜰리에구리", 3124482, "정인"),
    @(32, "건우씨엔에스", 222411, "정인"),
    @(33, "IDS코리아", 1892757, "정인"),
    @(34, "DB INC", 2886394, "정인"),
    @(35, "생생어르신복지센터", 2697039, "정인"),
    @(36, "알톤", 1896906, "민아"),
    @(37, "더드림핑", 1800688, "민아"),
    @(38, "한양패키지", 607620, "민아"),
    @(39, "펫츠비통", 4180383, "민아"),
    @(40, "해피발스데이", 3536711, "민아"),
    @(41, "휴앤고", 1281834, "민아"),
    @(42, "굿데이 남북결혼", 1836820, "민아"),
    @(43, "판타스틱코인노래방", 2310274, "민아"),
    @(44, "머크코리아 1-1", 1875332, "민아"),
    @(45, "머크코리아 2-1", 1254716, "민아"),
    @(46, "머크코리아 3-1", 1593239, "민아"),
    @(47, "머크코리아 GFA", 3671148, "민아"),
    @(48, "한우이츠", 3459602, "민아"),
    @(49, "한우이츠", 4005081, "민아"),
    @(50, "일두", 2901614, "민아"),
    @(51, "신한라이프케어", 1711939, "민아"),
    @(52, "SK텔레콤", 1725254, "민아")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Restore the selection state left behind in the sheet after the edit.
$ws.Range("D17:D18").Select()
